$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.410.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.382.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "407.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +10.30%  "
$ws.Range("E7").Value = "  +1.14%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.673"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.120"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "43.10"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.31%  "
$ws.Range("E12").Value = "  -1.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.903.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.70"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.400.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.356.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.31%  "
$ws.Range("E20").Value = "  -4.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "83.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "313.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.31%  "
$ws.Range("E25").Value = "  -0.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "29.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.22%  "
$ws.Range("E29").Value = "  -3.14%  "
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "41.25"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0480"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.95%  "
$ws.Range("E40").Value = "  -2.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "138.08"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.84%  "
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.299"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.124"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.09"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.57%  "
$ws.Range("E46").Value = "  -3.19%  "
$ws.Range("E47").Value = "  +2.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.130.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.33%  "
$ws.Range("E50").Value = "  -4.90%  "
$ws.Range("E51").Value = "  +0.02%  "
